$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = -1.722252655576526
$ws.Range("C24").Value = -0.1406863740731751
$ws.Range("D24").Value = -1.693277639921575
$ws.Range("E24").Value = -2.040686374073175
$ws.Range("F24").Value = 0.3060818441300339
$ws.Range("G24").Value = -2.082638238000783
$ws.Range("H24").Value = -0.8071502518150595
$ws.Range("I24").Value = 2.059313625926825
$ws.Range("J24").Value = 19.29561748113917
$ws.Range("K24").Value = -19.18321641127324
$ws.Range("B25").Value = 1.58156628150335
$ws.Range("C25").Value = 0.0289750156549502
$ws.Range("D25").Value = -0.3184337184966495
$ws.Range("E25").Value = 2.02833449970656
$ws.Range("F25").Value = -0.3603855824242574
$ws.Range("G25").Value = 0.9151024037614661
$ws.Range("H25").Value = 3.781566281503351
$ws.Range("I25").Value = 21.0178701367157
$ws.Range("J25").Value = -17.46096375569671
$ws.Range("K25").Value = -3.862991559831116
$ws.Range("B26").Value = -1.5525912658484
$ws.Range("C26").Value = -1.9
$ws.Range("D26").Value = 0.446768218203209
$ws.Range("E26").Value = -1.941951863927608
$ws.Range("F26").Value = -0.6664638777418844
$ws.Range("G26").Value = 2.2
$ws.Range("H26").Value = 19.43630385521235
$ws.Range("I26").Value = -19.04253003720006
$ws.Range("J26").Value = -5.444557841334467
$ws.Range("K26").Value = -2.749039957006937
$ws.Range("B27").Value = -0.3474087341515997
$ws.Range("C27").Value = 1.999359484051609
$ws.Range("D27").Value = -0.3893605980792076
$ws.Range("E27").Value = 0.8861273881065159
$ws.Range("F27").Value = 3.7525912658484
$ws.Range("G27").Value = 20.98889512106075
$ws.Range("H27").Value = -17.48993877135166
$ws.Range("I27").Value = -3.891966575486066
$ws.Range("J27").Value = -1.196448691158537
$ws.Range("K27").Value = 0.1128003469303707
$ws.Range("B28").Value = 2.346768218203209
$ws.Range("C28").Value = -0.04195186392760797
$ws.Range("D28").Value = 1.233536122258116
$ws.Range("E28").Value = 4.1
$ws.Range("F28").Value = 21.33630385521235
$ws.Range("G28").Value = -17.14253003720006
$ws.Range("H28").Value = -3.544557841334466
$ws.Range("I28").Value = -0.8490399570069369
$ws.Range("J28").Value = 0.4602090810819703
$ws.Range("K28").Value = 1.98156889852946
$ws.Range("B29").Value = -2.388720082130817
$ws.Range("C29").Value = -1.113232095945093
$ws.Range("D29").Value = 1.753231781796791
$ws.Range("E29").Value = 18.98953563700914
$ws.Range("F29").Value = -19.48929825540327
$ws.Range("G29").Value = -5.891326059537676
$ws.Range("H29").Value = -3.195808175210146
$ws.Range("I29").Value = -1.886559137121239
$ws.Range("J29").Value = -0.3651993196737491
$ws.Range("K29").Value = -6.130940961804114
$ws.Range("B30").Value = 1.275487986185724
$ws.Range("C30").Value = 4.141951863927608
$ws.Range("D30").Value = 21.37825571913996
$ws.Range("E30").Value = -17.10057817327245
$ws.Range("F30").Value = -3.502605977406859
$ws.Range("G30").Value = -0.8070880930793289
$ws.Range("H30").Value = 0.5021609450095783
$ws.Range("I30").Value = 2.023520762457068
$ws.Range("J30").Value = -3.742220879673297
$ws.Range("K30").Value = 3.165469857548658
$ws.Range("B31").Value = 2.866463877741884
$ws.Range("C31").Value = 20.10276773295423
$ws.Range("D31").Value = -18.37606615945818
$ws.Range("E31").Value = -4.778093963592582
$ws.Range("F31").Value = -2.082576079265053
$ws.Range("G31").Value = -0.7733270411761453
$ws.Range("H31").Value = 0.7480327762713443
$ws.Range("I31").Value = -5.017708865859021
$ws.Range("J31").Value = 1.889981871362934
$ws.Range("K31").Value = -0.5203221319907101
$ws.Range("B32").Value = 17.23630385521235
$ws.Range("C32").Value = -21.24253003720006
$ws.Range("D32").Value = -7.644557841334466
$ws.Range("E32").Value = -4.949039957006937
$ws.Range("F32").Value = -3.63979091891803
$ws.Range("G32").Value = -2.11843110147054
$ws.Range("H32").Value = -7.884172743600905
$ws.Range("I32").Value = -0.9764820063789501
$ws.Range("J32").Value = -3.386786009732595
$ws.Range("K32").Value = -5.082042055103904
$ws.Range("B33").Value = -38.47883389241241
$ws.Range("C33").Value = -24.88086169654682
$ws.Range("D33").Value = -22.18534381221929
$ws.Range("E33").Value = -20.87609477413038
$ws.Range("F33").Value = -19.35473495668289
$ws.Range("G33").Value = -25.12047659881326
$ws.Range("H33").Value = -18.2127858615913
$ws.Range("I33").Value = -20.62308986494494
$ws.Range("J33").Value = -22.31834591031625
$ws.Range("K33").Value = -19.37902846813061
$ws.Range("B34").Value = 13.59797219586559
$ws.Range("C34").Value = 16.29349008019312
$ws.Range("D34").Value = 17.60273911828203
$ws.Range("E34").Value = 19.12409893572952
$ws.Range("F34").Value = 13.35835729359916
$ws.Range("G34").Value = 20.26604803082111
$ws.Range("H34").Value = 17.85574402746747
$ws.Range("I34").Value = 16.16048798209616
$ws.Range("J34").Value = 19.0998054242818
$ws.Range("K34").Value = 17.72090440315328
$ws.Range("B35").Value = 2.695517884327529
$ws.Range("C35").Value = 4.004766922416437
$ws.Range("D35").Value = 5.526126739863926
$ws.Range("E35").Value = -0.2396149022664389
$ws.Range("F35").Value = 6.668075834955516
$ws.Range("G35").Value = 4.257771831601872
$ws.Range("H35").Value = 2.562515786230562
$ws.Range("I35").Value = 5.501833228416203
$ws.Range("J35").Value = 4.12293220728769
$ws.Range("K35").Value = 5.627651391037844
$ws.Range("B36").Value = 1.309249038088907
$ws.Range("C36").Value = 2.830608855536397
$ws.Range("D36").Value = -2.935132786593968
$ws.Range("E36").Value = 3.972557950627987
$ws.Range("F36").Value = 1.562253947274342
$ws.Range("G36").Value = -0.1330020980969671
$ws.Range("H36").Value = 2.806315344088674
$ws.Range("I36").Value = 1.42741432296016
$ws.Range("J36").Value = 2.932133506710315
$ws.Range("K36").Value = 2.605412008371104
$ws.Range("B37").Value = 1.52135981744749
$ws.Range("C37").Value = -4.244381824682876
$ws.Range("D37").Value = 2.66330891253908
$ws.Range("E37").Value = 0.2530049091854352
$ws.Range("F37").Value = -1.442251136185874
$ws.Range("G37").Value = 1.497066305999766
$ws.Range("H37").Value = 0.1181652848712531
$ws.Range("I37").Value = 1.622884468621408
$ws.Range("J37").Value = 1.296162970282197
$ws.Range("K37").Value = 2.134262027593551
$ws.Range("B38").Value = -5.765741642130365
$ws.Range("C38").Value = 1.14194909509159
$ws.Range("D38").Value = -1.268354908262054
$ws.Range("E38").Value = -2.963610953633364
$ws.Range("F38").Value = -0.0242935114477234
$ws.Range("G38").Value = -1.403194532576236
$ws.Range("H38").Value = 0.1015246511739181
$ws.Range("I38").Value = -0.225196847165293
$ws.Range("J38").Value = 0.6129022101460611
$ws.Range("K38").Value = -2.06679770341519
$ws.Range("B39").Value = 6.907690737221955
$ws.Range("C39").Value = 4.497386733868311
$ws.Range("D39").Value = 2.802130688497001
$ws.Range("E39").Value = 5.741448130682642
$ws.Range("F39").Value = 4.362547109554129
$ws.Range("G39").Value = 5.867266293304283
$ws.Range("H39").Value = 5.540544794965072
$ws.Range("I39").Value = 6.378643852276426
$ws.Range("J39").Value = 3.698943938715175
$ws.Range("K39").Value = 4.934366632778122
$ws.Range("B40").Value = -2.410304003353644
$ws.Range("C40").Value = -4.105560048724954
$ws.Range("D40").Value = -1.166242606539313
$ws.Range("E40").Value = -2.545143627667827
$ws.Range("F40").Value = -1.040424443917672
$ws.Range("G40").Value = -1.367145942256883
$ws.Range("H40").Value = -0.529046884945529
$ws.Range("I40").Value = -3.20874679850678
$ws.Range("J40").Value = -1.973324104443833
$ws.Range("K40").Value = -0.2490198157013881
$ws.Range("B41").Value = -1.695256045371309
$ws.Range("C41").Value = 1.244061396814331
$ws.Range("D41").Value = -0.1348396243141821
$ws.Range("E41").Value = 1.369879559435973
$ws.Range("F41").Value = 1.043158061096761
$ws.Range("G41").Value = 1.881257118408115
$ws.Range("H41").Value = -0.7984427951531354
$ws.Range("I41").Value = 0.4369798989098115
$ws.Range("J41").Value = 2.161284187652257
$ws.Range("K41").Value = 2.527729380265284
$ws.Range("B42").Value = 2.939317442185641
$ws.Range("C42").Value = 1.560416421057127
$ws.Range("D42").Value = 3.065135604807282
$ws.Range("E42").Value = 2.738414106468071
$ws.Range("F42").Value = 3.576513163779425
$ws.Range("G42").Value = 0.896813250218174
$ws.Range("H42").Value = 2.132235944281121
$ws.Range("I42").Value = 3.856540233023566
$ws.Range("J42").Value = 4.222985425636594
$ws.Range("K42").Value = -1.252462014395462
$ws.Range("B43").Value = -1.378901021128513
$ws.Range("C43").Value = 0.1258181626216415
$ws.Range("D43").Value = -0.2009033357175696
$ws.Range("E43").Value = 0.6371957215937845
$ws.Range("F43").Value = -2.042504191967466
$ws.Range("G43").Value = -0.8070814979045196
$ws.Range("H43").Value = 0.9172227908379254
$ws.Range("I43").Value = 1.283667983450953
$ws.Range("J43").Value = -4.191779456581102
$ws.Range("K43").Value = -0.8310663749237079
$ws.Range("B44").Value = 1.504719183750155
$ws.Range("C44").Value = 1.177997685410944
$ws.Range("D44").Value = 2.016096742722298
$ws.Range("E44").Value = -0.6636031708389534
$ws.Range("F44").Value = 0.5718195232239935
$ws.Range("G44").Value = 2.296123811966439
$ws.Range("H44").Value = 2.662569004579467
$ws.Range("I44").Value = -2.81287843545259
$ws.Range("J44").Value = 0.5478346462048052
$ws.Range("B45").Value = -0.326721498339211
$ws.Range("C45").Value = 0.511377558972143
$ws.Range("D45").Value = -2.168322354589108
$ws.Range("E45").Value = -0.932899660526161
$ws.Range("F45").Value = 0.7914046282162839
$ws.Range("G45").Value = 1.157849820829312
$ws.Range("H45").Value = -4.317597619202744
$ws.Range("I45").Value = -0.9568845375453494
$ws.Range("B46").Value = 0.838099057311354
$ws.Range("C46").Value = -1.841600856249897
$ws.Range("D46").Value = -0.60617816218695
$ws.Range("E46").Value = 1.118126126555495
$ws.Range("F46").Value = 1.484571319168523
$ws.Range("G46").Value = -3.990876120863533
$ws.Range("H46").Value = -0.6301630392061384
$ws.Range("B47").Value = -2.679699913561251
$ws.Range("C47").Value = -1.444277219498304
$ws.Range("D47").Value = 0.2800270692441409
$ws.Range("E47").Value = 0.646472261857169
$ws.Range("F47").Value = -4.828975178174887
$ws.Range("G47").Value = -1.468262096517492
$ws.Range("B48").Value = 1.235422694062947
$ws.Range("C48").Value = 2.959726982805392
$ws.Range("D48").Value = 3.32617217541842
$ws.Range("E48").Value = -2.149275264613636
$ws.Range("F48").Value = 1.211437817043759
$ws.Range("B49").Value = 1.724304288742445
$ws.Range("C49").Value = 2.090749481355473
$ws.Range("D49").Value = -3.384697958676583
$ws.Range("E49").Value = -0.02398487701918839
$ws.Range("B50").Value = 0.3664451926130281
$ws.Range("C50").Value = -5.109002247419028
$ws.Range("D50").Value = -1.748289165761633
$ws.Range("B51").Value = -5.475447440032056
$ws.Range("C51").Value = -2.114734358374661
$ws.Range("B52").Value = 3.360713081657395
